$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "EE00"
$ws.Cells.Item(2, 3).Value = "ME25"
$ws.Cells.Item(2, 4).Value = "ME30"
$ws.Cells.Item(2, 5).Value = "ME62"
$ws.Cells.Item(2, 8).Value = 3.4
$ws.Cells.Item(2, 9).Value = 3.5
$ws.Cells.Item(2, 10).Value = 3.0
$ws.Cells.Item(2, 11).Value = 3.2

$ws.Cells.Item(3, 2).Value = "EE01"
$ws.Cells.Item(3, 3).Value = "ME04"
$ws.Cells.Item(3, 4).Value = "ME15"
$ws.Cells.Item(3, 5).Value = "ME31"
$ws.Cells.Item(3, 8).Value = 2.5
$ws.Cells.Item(3, 9).Value = 3.4
$ws.Cells.Item(3, 10).Value = 3.5
$ws.Cells.Item(3, 11).Value = 3.1

$ws.Cells.Item(4, 2).Value = "EE02"
$ws.Cells.Item(4, 3).Value = "EE23"
$ws.Cells.Item(4, 4).Value = "ME32"
$ws.Cells.Item(4, 5).Value = "ME59"
$ws.Cells.Item(4, 8).Value = 3.6
$ws.Cells.Item(4, 9).Value = 3.7
$ws.Cells.Item(4, 10).Value = 2.2
$ws.Cells.Item(4, 11).Value = 3.9

$ws.Cells.Item(5, 2).Value = "EE03"
$ws.Cells.Item(5, 3).Value = "ME03"
$ws.Cells.Item(5, 4).Value = "ME33"
$ws.Cells.Item(5, 5).Value = "ME61"
$ws.Cells.Item(5, 8).Value = 3.7
$ws.Cells.Item(5, 9).Value = 3.3
$ws.Cells.Item(5, 10).Value = 3.3
$ws.Cells.Item(5, 11).Value = 2.1

$ws.Cells.Item(6, 2).Value = "EE04"
$ws.Cells.Item(6, 3).Value = "ME34"
$ws.Cells.Item(6, 4).Value = "ME54"
$ws.Cells.Item(6, 8).Value = 3.8
$ws.Cells.Item(6, 9).Value = 3.4
$ws.Cells.Item(6, 10).Value = 2.4

$ws.Cells.Item(7, 2).Value = "EE05"
$ws.Cells.Item(7, 3).Value = "ME05"
$ws.Cells.Item(7, 4).Value = "ME26"
$ws.Cells.Item(7, 5).Value = "ME35"
$ws.Cells.Item(7, 8).Value = 3.9
$ws.Cells.Item(7, 9).Value = 2.5
$ws.Cells.Item(7, 10).Value = 3.6
$ws.Cells.Item(7, 11).Value = 2.5

$ws.Cells.Item(8, 2).Value = "EE06"
$ws.Cells.Item(8, 3).Value = "ME06"
$ws.Cells.Item(8, 4).Value = "ME36"
$ws.Cells.Item(8, 5).Value = "ME56"
$ws.Cells.Item(8, 8).Value = 3.0
$ws.Cells.Item(8, 9).Value = 3.6
$ws.Cells.Item(8, 10).Value = 3.6
$ws.Cells.Item(8, 11).Value = 3.6

$ws.Cells.Item(9, 2).Value = "EE09"
$ws.Cells.Item(9, 3).Value = "ME07"
$ws.Cells.Item(9, 4).Value = "ME37"
$ws.Cells.Item(9, 5).Value = "ME58"
$ws.Cells.Item(9, 8).Value = 2.3
$ws.Cells.Item(9, 9).Value = 3.7
$ws.Cells.Item(9, 10).Value = 3.7
$ws.Cells.Item(9, 11).Value = 3.8

$ws.Cells.Item(10, 2).Value = "EE08"
$ws.Cells.Item(10, 3).Value = "ME02"
$ws.Cells.Item(10, 4).Value = "ME08"
$ws.Cells.Item(10, 5).Value = "ME38"
$ws.Cells.Item(10, 8).Value = 3.2
$ws.Cells.Item(10, 9).Value = 2.2
$ws.Cells.Item(10, 10).Value = 3.8
$ws.Cells.Item(10, 11).Value = 3.8

$ws.Cells.Item(11, 2).Value = "EE07"
$ws.Cells.Item(11, 3).Value = "ME09"
$ws.Cells.Item(11, 4).Value = "ME39"
$ws.Cells.Item(11, 5).Value = "ME60"
$ws.Cells.Item(11, 8).Value = 3.1
$ws.Cells.Item(11, 9).Value = 2.9
$ws.Cells.Item(11, 10).Value = 3.9
$ws.Cells.Item(11, 11).Value = 3.0

$ws.Cells.Item(12, 2).Value = "EE10"
$ws.Cells.Item(12, 3).Value = "ME10"
$ws.Cells.Item(12, 4).Value = "ME40"
$ws.Cells.Item(12, 5).Value = "ME53"
$ws.Cells.Item(12, 8).Value = 3.4
$ws.Cells.Item(12, 9).Value = 3.0
$ws.Cells.Item(12, 10).Value = 3.0
$ws.Cells.Item(12, 11).Value = 3.3

$ws.Cells.Item(13, 2).Value = "EE11"
$ws.Cells.Item(13, 3).Value = "ME11"
$ws.Cells.Item(13, 4).Value = "ME27"
$ws.Cells.Item(13, 5).Value = "ME41"
$ws.Cells.Item(13, 8).Value = 3.5
$ws.Cells.Item(13, 9).Value = 3.1
$ws.Cells.Item(13, 10).Value = 2.7
$ws.Cells.Item(13, 11).Value = 3.1

$ws.Cells.Item(14, 2).Value = "EE12"
$ws.Cells.Item(14, 3).Value = "ME12"
$ws.Cells.Item(14, 4).Value = "ME42"
$ws.Cells.Item(14, 5).Value = "ME63"
$ws.Cells.Item(14, 8).Value = 3.6
$ws.Cells.Item(14, 9).Value = 3.2
$ws.Cells.Item(14, 10).Value = 2.2
$ws.Cells.Item(14, 11).Value = 3.3

$ws.Cells.Item(15, 2).Value = "CpE00"
$ws.Cells.Item(15, 3).Value = "CpE01"
$ws.Cells.Item(15, 4).Value = "EE13"
$ws.Cells.Item(15, 5).Value = "ME43"
$ws.Cells.Item(15, 8).Value = 2.6
$ws.Cells.Item(15, 9).Value = 3.7
$ws.Cells.Item(15, 10).Value = 2.7
$ws.Cells.Item(15, 11).Value = 3.3

$ws.Cells.Item(16, 2).Value = "EE14"
$ws.Cells.Item(16, 3).Value = "ME14"
$ws.Cells.Item(16, 4).Value = "ME44"
$ws.Cells.Item(16, 5).Value = "ME57"
$ws.Cells.Item(16, 8).Value = 3.8
$ws.Cells.Item(16, 9).Value = 2.4
$ws.Cells.Item(16, 10).Value = 2.4
$ws.Cells.Item(16, 11).Value = 3.7

$ws.Cells.Item(17, 2).Value = "CpE02"
$ws.Cells.Item(17, 3).Value = "EE15"
$ws.Cells.Item(17, 4).Value = "ME01"
$ws.Cells.Item(17, 5).Value = "ME45"
$ws.Cells.Item(17, 8).Value = 3.8
$ws.Cells.Item(17, 9).Value = 3.9
$ws.Cells.Item(17, 10).Value = 2.1
$ws.Cells.Item(17, 11).Value = 3.5

$ws.Cells.Item(18, 2).Value = "EE16"
$ws.Cells.Item(18, 3).Value = "ME16"
$ws.Cells.Item(18, 4).Value = "ME24"
$ws.Cells.Item(18, 5).Value = "ME46"
$ws.Cells.Item(18, 8).Value = 3.0
$ws.Cells.Item(18, 9).Value = 3.6
$ws.Cells.Item(18, 10).Value = 3.4
$ws.Cells.Item(18, 11).Value = 3.6

$ws.Cells.Item(19, 2).Value = "EE17"
$ws.Cells.Item(19, 3).Value = "ME13"
$ws.Cells.Item(19, 4).Value = "ME17"
$ws.Cells.Item(19, 5).Value = "ME47"
$ws.Cells.Item(19, 6).Value = "ME55"
$ws.Cells.Item(19, 8).Value = 3.1
$ws.Cells.Item(19, 9).Value = 3.3
$ws.Cells.Item(19, 10).Value = 3.7
$ws.Cells.Item(19, 11).Value = 3.7
$ws.Cells.Item(19, 12).Value = 3.5

$ws.Cells.Item(20, 2).Value = "EE18"
$ws.Cells.Item(20, 3).Value = "ME00"
$ws.Cells.Item(20, 4).Value = "ME18"
$ws.Cells.Item(20, 5).Value = "ME48"
$ws.Cells.Item(20, 8).Value = 3.2
$ws.Cells.Item(20, 9).Value = 3.0
$ws.Cells.Item(20, 10).Value = 3.8
$ws.Cells.Item(20, 11).Value = 2.8

$ws.Cells.Item(21, 2).Value = "CpE03"
$ws.Cells.Item(21, 3).Value = "EE19"
$ws.Cells.Item(21, 4).Value = "ME19"
$ws.Cells.Item(21, 5).Value = "ME49"
$ws.Cells.Item(21, 8).Value = 3.9
$ws.Cells.Item(21, 9).Value = 2.3
$ws.Cells.Item(21, 10).Value = 3.9
$ws.Cells.Item(21, 11).Value = 3.9

$ws.Cells.Item(22, 2).Value = "EE20"
$ws.Cells.Item(22, 3).Value = "ME20"
$ws.Cells.Item(22, 4).Value = "ME23"
$ws.Cells.Item(22, 5).Value = "ME50"
$ws.Cells.Item(22, 8).Value = 3.4
$ws.Cells.Item(22, 9).Value = 3.0
$ws.Cells.Item(22, 10).Value = 3.3
$ws.Cells.Item(22, 11).Value = 3.0

$ws.Cells.Item(23, 2).Value = "EE21"
$ws.Cells.Item(23, 3).Value = "ME21"
$ws.Cells.Item(23, 4).Value = "ME29"
$ws.Cells.Item(23, 5).Value = "ME51"
$ws.Cells.Item(23, 8).Value = 3.5
$ws.Cells.Item(23, 9).Value = 2.1
$ws.Cells.Item(23, 10).Value = 3.9
$ws.Cells.Item(23, 11).Value = 3.1

$ws.Cells.Item(24, 2).Value = "EE22"
$ws.Cells.Item(24, 3).Value = "EE24"
$ws.Cells.Item(24, 4).Value = "ME22"
$ws.Cells.Item(24, 5).Value = "ME28"
$ws.Cells.Item(24, 6).Value = "ME52"
$ws.Cells.Item(24, 8).Value = 3.6
$ws.Cells.Item(24, 9).Value = 3.8
$ws.Cells.Item(24, 10).Value = 3.2
$ws.Cells.Item(24, 11).Value = 3.8
$ws.Cells.Item(24, 12).Value = 3.2

$ws.Cells.Item(25, 2).Value = "CE03"
$ws.Cells.Item(25, 3).Value = "CE04"
$ws.Cells.Item(25, 4).Value = "CE05"
$ws.Cells.Item(25, 5).Value = "CE06"
$ws.Cells.Item(25, 8).Value = 3.2
$ws.Cells.Item(25, 9).Value = 3.3
$ws.Cells.Item(25, 10).Value = 3.4
$ws.Cells.Item(25, 11).Value = 3.5

$ws.Cells.Item(26, 2).Value = "CE00"
$ws.Cells.Item(26, 3).Value = "CE01"
$ws.Cells.Item(26, 4).Value = "CE02"
$ws.Cells.Item(26, 8).Value = 3.9
$ws.Cells.Item(26, 9).Value = 2.7
$ws.Cells.Item(26, 10).Value = 3.1
